$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 31
$ws.Range("H31").Value = 1238.75
$ws.Range("I31").Value = 1258
$ws.Range("K31").Value = 3774
$ws.Range("M31").Value = -3544
# row 43
$ws.Range("H43").Value = 1751.2222
$ws.Range("J43").Value = 1480.5
$ws.Range("L43").Value = 1480.5
$ws.Range("N43").Value = -1618.5
# row 132
$ws.Range("H132").Value = 15153740
$ws.Range("I132").Value = 15386867
$ws.Range("K132").Value = 46160601
$ws.Range("M132").Value = -46158071
# row 135
$ws.Range("H135").Value = 1066.5143
$ws.Range("I135").Value = 691.25806
$ws.Range("K135").Value = 6221.32254
$ws.Range("M135").Value = -3686.32254
# row 137
$ws.Range("H137").Value = 5527.421
$ws.Range("I137").Value = 6482.8335
$ws.Range("J137").Value = 5086.4614
$ws.Range("K137").Value = 19448.5005
$ws.Range("L137").Value = 15259.3842
$ws.Range("M137").Value = -16898.5005
$ws.Range("N137").Value = -20359.3842
# row 141
$ws.Range("H141").Value = 12372.777
$ws.Range("I141").Value = 5774.619
$ws.Range("K141").Value = 17323.857
$ws.Range("M141").Value = -12143.857

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3397.341
$ws.Range("I32").Value = 2042.2958
$ws.Range("K32").Value = 2042.2958
$ws.Range("M32").Value = -1755.2958
# row 37
$ws.Range("H37").Value = 39495
$ws.Range("J37").Value = 44990
$ws.Range("L37").Value = 44990
$ws.Range("N37").Value = -45536
# row 61
$ws.Range("H61").Value = 3463.2856
$ws.Range("I61").Value = 3057.5
$ws.Range("J61").Value = 4004.3333
$ws.Range("K61").Value = 3057.5
$ws.Range("L61").Value = 4004.3333
$ws.Range("M61").Value = -2845.5
$ws.Range("N61").Value = -4428.3333
# row 63
$ws.Range("H63").Value = 2447.0625
$ws.Range("I63").Value = 2375.2144
$ws.Range("K63").Value = 2375.2144
$ws.Range("M63").Value = -1689.2144
# row 66
$ws.Range("H66").Value = 2447.0625
$ws.Range("I66").Value = 2375.2144
$ws.Range("K66").Value = 11876.072
$ws.Range("M66").Value = -8444.072
# row 74
$ws.Range("H74").Value = 64958.855
$ws.Range("I74").Value = 1197.2222
$ws.Range("K74").Value = 1197.2222
$ws.Range("M74").Value = -323.2221999999999
# row 77
$ws.Range("H77").Value = 64958.855
$ws.Range("I77").Value = 1197.2222
$ws.Range("K77").Value = 5986.111
$ws.Range("M77").Value = -1618.111
# row 97
$ws.Range("H97").Value = 1712884.1
$ws.Range("I97").Value = 1802488.2
$ws.Range("K97").Value = 1802488.2
$ws.Range("M97").Value = -1801992.2
# row 132
$ws.Range("H132").Value = 2935.7856
$ws.Range("J132").Value = 4497.5
$ws.Range("L132").Value = 13492.5
$ws.Range("N132").Value = -18552.5
# row 136
$ws.Range("H136").Value = 3463.2856
$ws.Range("I136").Value = 3057.5
$ws.Range("J136").Value = 4004.3333
$ws.Range("K136").Value = 9172.5
$ws.Range("L136").Value = 12012.9999
$ws.Range("M136").Value = -6622.5
$ws.Range("N136").Value = -17112.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 1256.625
$ws.Range("I20").Value = 1165
$ws.Range("K20").Value = 1165
$ws.Range("M20").Value = -918
# row 107
$ws.Range("H107").Value = 2383045
$ws.Range("I107").Value = 3403459.2
$ws.Range("K107").Value = 3403459.2
$ws.Range("M107").Value = -3401539.2
# row 134
$ws.Range("H134").Value = 1837.0182
$ws.Range("I134").Value = 789.25
$ws.Range("J134").Value = 19998.334
$ws.Range("K134").Value = 2367.75
$ws.Range("L134").Value = 59995.00199999999
$ws.Range("M134").Value = 167.25
$ws.Range("N134").Value = -65065.00199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 291.125
$ws.Range("I7").Value = 193.875
$ws.Range("K7").Value = 193.875
$ws.Range("M7").Value = -80.875
# row 58
$ws.Range("H58").Value = 1744.8182
$ws.Range("I58").Value = 1590.0476
$ws.Range("K58").Value = 1590.0476
$ws.Range("M58").Value = -1387.0476
# row 132
$ws.Range("H132").Value = 2777.0645
$ws.Range("I132").Value = 2608.125
$ws.Range("J132").Value = 3356.2856
$ws.Range("K132").Value = 7824.375
$ws.Range("L132").Value = 10068.8568
$ws.Range("M132").Value = -5294.375
$ws.Range("N132").Value = -15128.8568
# row 134
$ws.Range("H134").Value = 2994.1072
$ws.Range("I134").Value = 2633.6
$ws.Range("K134").Value = 7900.799999999999
$ws.Range("M134").Value = -5365.799999999999
# row 135
$ws.Range("H135").Value = 93874.25
$ws.Range("J135").Value = 93874.25
$ws.Range("L135").Value = 93874.25
$ws.Range("N135").Value = -104014.25
# row 136
$ws.Range("H136").Value = 1744.8182
$ws.Range("I136").Value = 1590.0476
$ws.Range("K136").Value = 4770.142800000001
$ws.Range("M136").Value = -2220.142800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 132
$ws.Range("H132").Value = 2175.074
$ws.Range("J132").Value = 2468.182
$ws.Range("L132").Value = 22213.638
$ws.Range("N132").Value = -27273.638
# row 134
$ws.Range("H134").Value = 1622.4
$ws.Range("I134").Value = 1622.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4867.200000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 202.7999999999993
$ws.Range("N134").ClearContents()
# row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
# row 138
$ws.Range("H138").Value = 13829.235
$ws.Range("I138").Value = 13829.235
$ws.Range("K138").Value = 41487.705
$ws.Range("M138").Value = -36347.705

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 851313.7
$ws.Range("I97").Value = 916717.75
$ws.Range("K97").Value = 916717.75
$ws.Range("M97").Value = -916221.75
# row 123
$ws.Range("H123").Value = 35999
$ws.Range("J123").Value = 35999
$ws.Range("L123").Value = 35999
$ws.Range("N123").Value = -40899

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 3720.7896
$ws.Range("I68").Value = 3545.4546
$ws.Range("K68").Value = 3545.4546
$ws.Range("M68").Value = -2796.4546
# row 71
$ws.Range("H71").Value = 3720.7896
$ws.Range("I71").Value = 3545.4546
$ws.Range("K71").Value = 17727.273
$ws.Range("M71").Value = -13983.273
# row 93
$ws.Range("H93").Value = 66673068
$ws.Range("I93").Value = 66673068
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 66673068
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -66671820
$ws.Range("N93").ClearContents()
# row 100
$ws.Range("H100").Value = 4984.3076
$ws.Range("I100").Value = 4984.3076
$ws.Range("K100").Value = 4984.3076
$ws.Range("M100").Value = -4443.3076
# row 122
$ws.Range("H122").Value = 4929.2573
$ws.Range("J122").Value = 7485.4443
$ws.Range("L122").Value = 22456.3329
$ws.Range("N122").Value = -27356.3329
# row 132
$ws.Range("H132").Value = 2913.8035
$ws.Range("I132").Value = 2303.58
$ws.Range("K132").Value = 6910.74
$ws.Range("M132").Value = -4380.74

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 46
$ws.Range("H46").Value = 66666
$ws.Range("J46").Value = 66666
$ws.Range("L46").Value = 66666
$ws.Range("N46").Value = -67128
# row 122
$ws.Range("H122").Value = 3543.1667
$ws.Range("I122").Value = 2759.1875
$ws.Range("K122").Value = 8277.5625
$ws.Range("M122").Value = -5827.5625
# row 132
$ws.Range("H132").Value = 19608856
$ws.Range("I132").Value = 21739974
$ws.Range("K132").Value = 65219922
$ws.Range("M132").Value = -65217392
# row 134
$ws.Range("H134").Value = 66666
$ws.Range("J134").Value = 66666
$ws.Range("L134").Value = 199998
$ws.Range("N134").Value = -205068
# row 135
$ws.Range("H135").Value = 103106.305
$ws.Range("J135").Value = 104638.5
$ws.Range("L135").Value = 104638.5
$ws.Range("N135").Value = -114778.5
# row 136
$ws.Range("H136").Value = 1535.2821
$ws.Range("I136").Value = 1060.9333
$ws.Range("J136").Value = 3116.4443
$ws.Range("K136").Value = 3182.7999
$ws.Range("L136").Value = 9349.332900000001
$ws.Range("M136").Value = -632.7999
$ws.Range("N136").Value = -14449.3329
# row 141
$ws.Range("H141").Value = 118247.25
$ws.Range("J141").Value = 118247.25
$ws.Range("L141").Value = 118247.25
$ws.Range("N141").Value = -128607.25
